$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 293, shifting existing rows 293-375 down to 295-377
$ws.Rows("293:294").Insert()

# Row 293: new record
$ws.Range("A293").Value = 4
$ws.Range("B293").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C293").Value = 'Los Lagos'
$ws.Range("D293").Value = 44964
$ws.Range("E293").Value = 10
$ws.Range("F293").Value = 'Fruta'
$ws.Range("G293").Value = 100101
$ws.Range("H293").Value = 'Berries'
$ws.Range("I293").Value = 100112025
$ws.Range("J293").Value = 'Frutilla'
$ws.Range("K293").Value = 'Sin especificar'
$ws.Range("L293").Value = 'Primera'
$ws.Range("M293").Value = 600
$ws.Range("N293").Value = 9000
$ws.Range("O293").Value = 10000
$ws.Range("P293").Value = 9500
$ws.Range("Q293").Value = '$/caja 7 kilos'
$ws.Range("R293").Value = 'Región de La Araucanía'
$ws.Range("S293").Value = 1357
$ws.Range("T293").Value = 7

# Row 294: new record
$ws.Range("A294").Value = 4
$ws.Range("B294").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C294").Value = 'Los Lagos'
$ws.Range("D294").Value = 44964
$ws.Range("E294").Value = 10
$ws.Range("F294").Value = 'Fruta'
$ws.Range("G294").Value = 100101
$ws.Range("H294").Value = 'Berries'
$ws.Range("I294").Value = 100112025
$ws.Range("J294").Value = 'Frutilla'
$ws.Range("K294").Value = 'Sin especificar'
$ws.Range("L294").Value = 'Segunda'
$ws.Range("M294").Value = 200
$ws.Range("N294").Value = 8000
$ws.Range("O294").Value = 8000
$ws.Range("P294").Value = 8000
$ws.Range("Q294").Value = '$/caja 7 kilos'
$ws.Range("R294").Value = 'Región de La Araucanía'
$ws.Range("S294").Value = 1143
$ws.Range("T294").Value = 7
